# CSCI_1108_Lab01C_Spirograph.docx edit script
# "drawing a square and a spirograph"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Tidy the "go3 button draws" sentence: Word's grammar checker had
#    flagged "draws" (gramStart/gramEnd); clearing that suggestion also
#    collapses the three runs that made up the sentence into one.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "blue pattern and the go3 button draws the orange pattern.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "blue pattern and the go3 button draws the orange pattern.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Clear the grammar flag around "spirograph.nlogo".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "spirograph.nlogo",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "spirograph.nlogo", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Merge the "1" / "3" runs that spell out "13" (points value).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "B: 13 points]:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "B: 13 points]:", 2) | Out-Null

Write-Output "done"
